$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to rounded (custom accuracy) figures
$ws.Range("B5").Value = 3.84
$ws.Range("C5").Value = 2.55
$ws.Range("D5").Value = 0.69
$ws.Range("E5").Value = 8.05
$ws.Range("F5").Value = 6.66
$ws.Range("G5").Value = 3.03
$ws.Range("H5").Value = 14.64
$ws.Range("I5").Value = 4.65
$ws.Range("J5").Value = 1.96
$ws.Range("K5").Value = 2.91
$ws.Range("L5").Value = 3.33
$ws.Range("M5").Value = 3.34
$ws.Range("N5").Value = 0.97
$ws.Range("O5").Value = 3.01
$ws.Range("P5").Value = 4.21
$ws.Range("Q5").Value = 2.73
$ws.Range("R5").Value = 0.73
$ws.Range("S5").Value = 0.34
$ws.Range("T5").Value = 38.48
$ws.Range("U5").Value = 8.57
$ws.Range("V5").Value = 2.78
$ws.Range("W5").Value = 5.57
$ws.Range("X5").Value = 3.08
$ws.Range("Y5").Value = 0.29
$ws.Range("Z5").Value = 6.66
$ws.Range("AA5").Value = 2.45
$ws.Range("AB5").Value = 2.3
$ws.Range("AC5").Value = 2.68
$ws.Range("AD5").Value = 3.42
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 13.3
$ws.Range("AG5").Value = 1.47
$ws.Range("AH5").Value = 3.47

# Remove the now-extra data row (row 6), shifting dimension from AH6 to AH5
$ws.Rows.Item(6).Delete()
